# Update the "想去人数" (interested-in-going count) figures in the
# "展览" and "全部类型" sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 14
    $ws.Range("F9").Value = 3885
    $ws.Range("F10").Value = 4299
    $ws.Range("F12").Value = 130
}
